$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "En propiedad, por herencia o donación"
$ws.Range("A2").Value = "En alquiler"
$ws.Range("A3").Value = "Otra forma"
$ws.Range("A4").Value = "En propiedad, por compra totalmente pagada"
$ws.Range("A5").Value = "Cedida gratis o a bajo precio por otro hogar, la empresa,¿"
$ws.Range("A6").Value = "En propiedad, con pagos pendientes (hipotecas),¿"
